$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New quarterly row appended below the existing data (row 63).
# Force column A to be treated as text first so the date-like string
# "01-04-2021" is stored verbatim (as a shared string) instead of being
# auto-converted into a date serial number, then reset the style so the
# cell carries no explicit style index (matching the rest of column A).
$ws.Range("A63").NumberFormat = "@"
$ws.Range("A63").Value = "01-04-2021"
$ws.Range("A63").Style = "Normal"

$ws.Range("B63").Value = 33.1
$ws.Range("C63").Value = 15.6
$ws.Range("D63").Value = 18.4
$ws.Range("E63").Value = -2.9
$ws.Range("F63").Value = 51.3
$ws.Range("G63").Value = 12.7
